$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Free up the shared-string slot held by "Afghanistan" (C5) before introducing
# the new "ExtraActionIndicator" header text, so the new header text reuses
# that freed slot (matches the original authoring order).
$ws.Range("C5").Value = ""

# New column D header: "ExtraActionIndicator"
$ws.Range("D1").Value = "ExtraActionIndicator"

# Fill column D down to row 4 with "NULL"
$ws.Range("D2").Value = """NULL"""
$ws.Range("D3").Value = """NULL"""
$ws.Range("D4").Value = """NULL"""

# Row 5: country value changes from "Afghanistan" to "Bangladesh",
# and the new ExtraActionIndicator column gets "INSTANCENO:2"
$ws.Range("C5").Value = """Bangladesh"""
$ws.Range("D5").Value = """INSTANCENO:2"""

# Remaining rows in column D get "NULL"
$ws.Range("D6").Value = """NULL"""
$ws.Range("D7").Value = """NULL"""
$ws.Range("D8").Value = """NULL"""
$ws.Range("D9").Value = """NULL"""
$ws.Range("D10").Value = """NULL"""
$ws.Range("D11").Value = """NULL"""
$ws.Range("D12").Value = """NULL"""

# Give the new header cell the same bold font as the other header cells
# (A1:C1), which creates a distinct style entry analogous to the existing
# bold-with-alignment-flag style already present in the workbook.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").WrapText = $false

# Column D width, close to the recorded best-fit width
$ws.Columns.Item(4).ColumnWidth = 18.5

# Update the active selection to D12, matching the new end of the data range
$ws.Range("D12").Select()
